$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $newValue
    $c.Style = "Normal"
}

Set-TextValue "D2" "243.73"
Set-TextValue "D3" "23.02"
Set-TextValue "D4" "5.400"
Set-TextValue "D6" "3.432"
Set-TextValue "D7" "6.494"
Set-TextValue "D8" "0.8110"
Set-TextValue "D9" "0.9270"
Set-TextValue "D11" "0.07386"
Set-TextValue "D12" "0.03177"
Set-TextValue "D13" "0.03080"
Set-TextValue "D14" "0.09349"
Set-TextValue "D16" "0.001576"
Set-TextValue "D17" "0.04701"
Set-TextValue "B18" "TigerCash"
Set-TextValue "C18" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D18" "0.005969"
Set-TextValue "E18" "17TigerCashTCH"
Set-TextValue "B19" "BitKan"
Set-TextValue "C19" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D19" "0.001257"
Set-TextValue "E19" "18BitKanKAN"
Set-TextValue "B20" "HotbitToken"
Set-TextValue "C20" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D20" "0.004789"
Set-TextValue "E20" "19HotbitTokenHTB"
Set-TextValue "B21" "NitroEx"
Set-TextValue "C21" "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D21" "0.00008005"
Set-TextValue "E21" "20NitroExNTX"
Set-TextValue "B22" "LEO"
Set-TextValue "C22" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D22" "3.557"
Set-TextValue "E22" "21LEOLEO"
Set-TextValue "B23" "BTSEToken"
Set-TextValue "C23" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D23" "2.133"
Set-TextValue "E23" "22BTSETokenBTSE"
Set-TextValue "B24" "One"
Set-TextValue "C24" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D24" "0.01113"
Set-TextValue "E24" "23OneONEBestin24h"
Set-TextValue "D26" "0.1331"
Set-TextValue "D27" "0.0002341"
Set-TextValue "E27" "26UpBotsUBXTWorstin24h"
Set-TextValue "D40" "0.03924"
Set-TextValue "D41" "0.006334"
Set-TextValue "E41" "40KickTokenKICK"
Set-TextValue "D43" "0.003502"
Set-TextValue "D44" "0.008363"
Set-TextValue "D45" "0.00005195"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "D47" "0.6664"
Set-TextValue "D48" "0.002064"
Set-TextValue "D49" "0.00002101"
Set-TextValue "D50" "0.0002001"
